## Reapply "new changes in ops (ordercreation & orderpage & order form)"
## Restructures the order-tracking sheet: reorders/renames a few columns,
## adds a new "Tier" column (N), and replaces the sample order rows with
## new sample data (new OrderIDs, Emp IDs, State/County, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Style plumbing first (done before values so the run's shared-string
#    table ends up holding only the strings that survive into the final
#    grid -- old labels that disappear entirely, e.g. "SIPL0004" /
#    "SIPL0005", are simply never written again).
# ---------------------------------------------------------------------

# C2/C3, G2/G3 and D2 all need the "font19, no fill, bordered" look that
# E2/E3 already carry in the source sheet -- clone it.
$ws.Range("E2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)

# D2 needs a border with no left edge (new borderId) -- build it by hand
# on top of the font19 look just applied, then stamp it onto D3.
$ws.Range("D2").Borders.Item(7).LineStyle = -4142
$ws.Range("D2").Borders.Item(8).LineStyle = 1
$ws.Range("D2").Borders.Item(8).ColorIndex = 1
$ws.Range("D2").Borders.Item(9).LineStyle = 1
$ws.Range("D2").Borders.Item(9).ColorIndex = 1
$ws.Range("D2").Borders.Item(10).LineStyle = 1
$ws.Range("D2").Borders.Item(10).ColorIndex = 1
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# E2/E3 drop back down to the plain body look (style formerly on B2).
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)

# New column N: header cell matches the other header cells, body cells
# match the plain body look.
$ws.Range("A1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Header row (row 1) -- same seven columns keep the same labels; the
#    Typist/Typist QC/Client/Lob/Product Name group gets reordered, and
#    "Tier" is appended as a new column N.
# ---------------------------------------------------------------------

$ws.Range("E1").Value2 = "Typist"
$ws.Range("F1").Value2 = "Typist QC"
$ws.Range("G1").Value2 = "Client"
$ws.Range("H1").Value2 = "Lob"
$ws.Range("I1").Value2 = "Process"
$ws.Range("J1").Value2 = "Product Name"
$ws.Range("N1").Value2 = "Tier"

# ---------------------------------------------------------------------
# 3) Row 2 data
# ---------------------------------------------------------------------

$ws.Range("B2").Value2 = "FINN18-001"
$ws.Range("C2").Value2 = "SIPL4167"
$ws.Range("D2").Value2 = "SIPL5688"
$ws.Range("E2").Value2 = "SIPL5317"
$ws.Range("F2").Value2 = "SIPL5317"
$ws.Range("G2").Value2 = "FINN TITLE"
$ws.Range("H2").Value2 = "Title"
$ws.Range("I2").Value2 = "Search & Typing"
$ws.Range("J2").Value2 = "Property Reports"
$ws.Range("K2").Value2 = "AL"
$ws.Range("L2").Value2 = "Shelby"
$ws.Range("M2").Value2 = "WIP"

# ---------------------------------------------------------------------
# 4) Row 3 data
# ---------------------------------------------------------------------

$ws.Range("B3").Value2 = "FINN18-002"
$ws.Range("C3").Value2 = "SIPL6153"
$ws.Range("D3").Value2 = "SIPL5688"
$ws.Range("E3").Value2 = "SIPL0102"
$ws.Range("F3").Value2 = "SIPL0103"
$ws.Range("G3").Value2 = "FINN TITLE"
$ws.Range("H3").Value2 = "Title"
$ws.Range("I3").Value2 = "Search & Typing"
$ws.Range("J3").Value2 = "Foreclosure information Report"
$ws.Range("K3").Value2 = "FL"
$ws.Range("L3").Value2 = "Clay"
$ws.Range("M3").Value2 = "WIP"

# ---------------------------------------------------------------------
# 5) Column widths -- column C shrinks (less sprawling content), the
#    "State"-ish mid columns shift over a couple of slots, a wide new
#    column J needs room for long product names, and new column N gets
#    a modest width.
# ---------------------------------------------------------------------

$ws.Columns("C").ColumnWidth = 19.6
$ws.Columns("G").ColumnWidth = 11.83
$ws.Columns("H").ColumnWidth = 15.1
$ws.Columns("I").ColumnWidth = 15.1
$ws.Columns("J").ColumnWidth = 28.2
$ws.Columns("N").ColumnWidth = 10.7

# ---------------------------------------------------------------------
# 6) Selection cosmetics to match the saved view.
# ---------------------------------------------------------------------

$ws.Range("H6").Select() | Out-Null
